$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 (shifts rows 7-20 down to 8-21)
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 347
$ws.Range("B7").Value = "Zelda Remastered"
$ws.Range("C7").Value = "Nintendo Switch"
$ws.Range("D7").Value = 500
$ws.Range("E7").Value = 1
